$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save the original values of the rows that are being cyclically shifted
# (row order before edit: 2, 3, 4, 6 ; row 5 is untouched)
$row2 = @{ D=44175; K="Rainier"; L="Segunda"; M=270;  N=25000; O=26000; P=25500; Q="$/caja 18 kilos";    R="Región de O'Higgins"; S=1417; T=18 }
$row3 = @{ D=44208; K="Lapins";  L="Segunda"; M=200;  N=10500; O=11000; P=10750; Q="$/bandeja 12 kilos"; R="Provincia de Curicó";  S=896;  T=12 }
$row4 = @{ D=44229; K="Santina"; L="Primera"; M=250;  N=6500;  O=7000;  P=6750;  Q="$/bandeja 5 kilos";  R="Provincia de Curicó";  S=1350; T=5  }
$row6 = @{ D=44161; K="Bing";    L="Primera"; M=160;  N=39000; O=40000; P=39500; Q="$/caja 20 kilos";    R="Provincia de Curicó";  S=1975; T=20 }

function Set-DataRow($rowNum, $data) {
    $ws.Range("D$rowNum").Value = $data.D
    $ws.Range("K$rowNum").Value = $data.K
    $ws.Range("L$rowNum").Value = $data.L
    $ws.Range("M$rowNum").Value = $data.M
    $ws.Range("N$rowNum").Value = $data.N
    $ws.Range("O$rowNum").Value = $data.O
    $ws.Range("P$rowNum").Value = $data.P
    $ws.Range("Q$rowNum").Value = $data.Q
    $ws.Range("R$rowNum").Value = $data.R
    $ws.Range("S$rowNum").Value = $data.S
    $ws.Range("T$rowNum").Value = $data.T
}

# Cyclic shift: new row2 <- old row6, new row3 <- old row2, new row4 <- old row3, new row6 <- old row4
Set-DataRow 2 $row6
Set-DataRow 3 $row2
Set-DataRow 4 $row3
Set-DataRow 6 $row4
